$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write the 17 "F:V" columns of a data row (home..url) in one shot.
function Set-RowFV($row, $vals) {
    $arr = New-Object 'object[,]' 1,17
    for ($i = 0; $i -lt 17; $i++) {
        $arr[0,$i] = $vals[$i]
    }
    $rng = $ws.Range("F" + $row + ":V" + $row)
    $rng.Value = $arr
}

# Helper: append a brand-new match row (A:V), copying number/style formatting
# from the last existing data row (181) first, then filling in the values.
function Add-NewRow($row, $indice, $eSerial, $vals) {
    $srcRow = $ws.Range("A181:V181")
    $dstRow = $ws.Range("A" + $row + ":V" + $row)
    $srcRow.Copy()
    $dstRow.PasteSpecial(-4122)

    $ws.Range("A" + $row).Value = $indice
    $ws.Range("B" + $row).Value = "italy"
    $ws.Range("C" + $row).Value = "serie-c-group-c"
    $ws.Range("D" + $row).Value = "2023-2024"
    $ws.Range("E" + $row).Value = $eSerial

    Set-RowFV $row $vals
}

# --- Re-order rows 95-181: swapped / cycled pairs of matches (same kick-off
#     slot, re-sorted) coming from the upstream scrape re-run. ---
Set-RowFV 95 @('Monopoli', 1, 'Picerno', 1, 2.59, '24/10/2023 12:42', 2.67, '25/10/2023 20:41', 2.92, '24/10/2023 12:42', 2.97, '25/10/2023 20:38', 2.69, '24/10/2023 12:42', 2.9, '25/10/2023 20:41', 'https://www.betexplorer.com/football/italy/serie-c-group-c/monopoli-picerno/jD1hGTen/')
Set-RowFV 96 @('Casertana', 2, 'Juve Stabia', 1, 2.56, '24/10/2023 12:42', 2.69, '25/10/2023 20:37', 2.83, '24/10/2023 12:42', 3.01, '25/10/2023 20:37', 2.8, '24/10/2023 12:42', 2.83, '25/10/2023 20:37', 'https://www.betexplorer.com/football/italy/serie-c-group-c/casertana-juve-stabia/vL2bC8UL/')
Set-RowFV 103 @('Audace Cerignola', 2, 'Casertana', 4, 1.72, '26/10/2023 09:12', 2.19, '29/10/2023 18:21', 3.22, '26/10/2023 09:12', 2.96, '29/10/2023 18:21', 4.74, '26/10/2023 09:12', 3.8, '29/10/2023 18:21', 'https://www.betexplorer.com/football/italy/serie-c-group-c/audace-cerignola-casertana/je62BSqS/')
Set-RowFV 104 @('Juve Stabia', 1, 'Latina', 0, 1.99, '27/10/2023 02:42', 2.01, '29/10/2023 18:26', 2.97, '27/10/2023 02:42', 3.28, '29/10/2023 18:26', 3.79, '27/10/2023 02:42', 3.92, '29/10/2023 18:21', 'https://www.betexplorer.com/football/italy/serie-c-group-c/juve-stabia-latina/pzwl53mo/')
Set-RowFV 105 @('Picerno', 2, 'Foggia', 0, 1.98, '27/10/2023 02:42', 3.15, '29/10/2023 18:26', 2.98, '27/10/2023 02:42', 2.75, '29/10/2023 18:26', 3.82, '27/10/2023 02:42', 2.66, '29/10/2023 18:26', 'https://www.betexplorer.com/football/italy/serie-c-group-c/picerno-foggia/nuHWTnQo/')
Set-RowFV 106 @('Turris', 0, 'Giugliano', 1, 2.13, '27/10/2023 02:42', 1.99, '29/10/2023 18:24', 3.07, '27/10/2023 02:42', 3.54, '29/10/2023 18:24', 3.27, '27/10/2023 02:42', 3.68, '29/10/2023 18:24', 'https://www.betexplorer.com/football/italy/serie-c-group-c/turris-giugliano/K6UwSQfb/')
Set-RowFV 107 @('Virtus Francavilla', 1, 'Taranto', 2, 2.15, '27/10/2023 02:42', 2.78, '29/10/2023 20:41', 2.95, '27/10/2023 02:42', 2.82, '29/10/2023 20:41', 3.36, '27/10/2023 02:42', 2.91, '29/10/2023 20:41', 'https://www.betexplorer.com/football/italy/serie-c-group-c/virtus-francavilla-taranto/bcTsRpA4/')
Set-RowFV 108 @('Brindisi', 1, 'Monopoli', 2, 2.15, '27/10/2023 02:42', 2.32, '29/10/2023 20:41', 2.95, '27/10/2023 02:42', 3.04, '29/10/2023 20:41', 3.36, '27/10/2023 02:42', 3.36, '29/10/2023 20:41', 'https://www.betexplorer.com/football/italy/serie-c-group-c/brindisi-monopoli/z1mEB5AH/')
Set-RowFV 117 @('Taranto', 0, 'Juve Stabia', 2, 2.83, '02/11/2023 08:12', 2.87, '05/11/2023 18:26', 2.8, '02/11/2023 08:12', 2.54, '05/11/2023 18:26', 2.57, '02/11/2023 08:12', 3.17, '05/11/2023 18:26', 'https://www.betexplorer.com/football/italy/serie-c-group-c/taranto-juve-stabia/fPrMv4v4/')
Set-RowFV 119 @('ACR Messina', 0, 'Benevento', 1, 3.65, '02/11/2023 08:12', 4.09, '05/11/2023 18:22', 3.08, '02/11/2023 08:12', 3.23, '05/11/2023 18:28', 1.99, '02/11/2023 08:12', 1.98, '05/11/2023 18:28', 'https://www.betexplorer.com/football/italy/serie-c-group-c/acr-messina-benevento/tMInQ4PA/')
Set-RowFV 127 @('Audace Cerignola', 1, 'Catania', 0, 2.64, '09/11/2023 09:13', 2.62, '12/11/2023 18:09', 2.8, '09/11/2023 09:13', 2.98, '12/11/2023 18:09', 2.74, '09/11/2023 09:13', 2.93, '12/11/2023 18:09', 'https://www.betexplorer.com/football/italy/serie-c-group-c/audace-cerignola-catania/2knUxr9G/')
Set-RowFV 128 @('Brindisi', 0, 'Avellino', 4, 4.7, '09/11/2023 09:13', 5.05, '12/11/2023 18:20', 3.22, '09/11/2023 09:13', 3.86, '12/11/2023 18:20', 1.73, '09/11/2023 09:13', 1.67, '12/11/2023 18:20', 'https://www.betexplorer.com/football/italy/serie-c-group-c/brindisi-avellino/Mwa2KZhk/')
Set-RowFV 137 @('Foggia', 2, 'ACR Messina', 0, 1.74, '16/11/2023 09:12', 1.66, '19/11/2023 15:50', 3.27, '16/11/2023 09:12', 3.47, '19/11/2023 16:03', 4.56, '16/11/2023 09:12', 5.35, '19/11/2023 15:50', 'https://www.betexplorer.com/football/italy/serie-c-group-c/foggia-acr-messina/2DYRIsZs/')
Set-RowFV 138 @('Avellino', 1, 'Giugliano', 3, 1.4, '16/11/2023 09:12', 1.3, '19/11/2023 16:11', 4.02, '16/11/2023 09:12', 5.21, '19/11/2023 16:11', 7.15, '16/11/2023 09:12', 10.42, '19/11/2023 16:11', 'https://www.betexplorer.com/football/italy/serie-c-group-c/avellino-giugliano/rTieO3BK/')
Set-RowFV 139 @('Monopoli', 3, 'Benevento', 0, 3.58, '16/11/2023 09:12', 3.22, '19/11/2023 16:10', 3.09, '16/11/2023 09:12', 3.14, '19/11/2023 16:14', 2.01, '16/11/2023 09:12', 2.33, '19/11/2023 16:11', 'https://www.betexplorer.com/football/italy/serie-c-group-c/monopoli-benevento/lUzvGuJ0/')
Set-RowFV 140 @('Potenza', 2, 'Audace Cerignola', 2, 2.29, '16/11/2023 09:12', 2.72, '19/11/2023 18:26', 2.94, '16/11/2023 09:12', 3.11, '19/11/2023 18:26', 3.07, '16/11/2023 09:12', 2.71, '19/11/2023 18:26', 'https://www.betexplorer.com/football/italy/serie-c-group-c/potenza-audace-cerignola/UwqnEJmD/')
Set-RowFV 141 @('Juve Stabia', 0, 'Sorrento', 0, 1.61, '16/11/2023 09:12', 1.5, '19/11/2023 18:28', 3.49, '16/11/2023 09:12', 3.78, '19/11/2023 18:28', 5.18, '16/11/2023 09:12', 8.04, '19/11/2023 18:28', 'https://www.betexplorer.com/football/italy/serie-c-group-c/juve-stabia-sorrento/t4XVH1ll/')
Set-RowFV 155 @('Audace Cerignola', 0, 'Picerno', 1, 2.32, '30/11/2023 09:13', 2.44, '03/12/2023 20:43', 2.91, '30/11/2023 09:13', 3.03, '03/12/2023 20:43', 3.13, '30/11/2023 09:13', 3.14, '03/12/2023 20:43', 'https://www.betexplorer.com/football/italy/serie-c-group-c/audace-cerignola-picerno/bJtCktAJ/')
Set-RowFV 156 @('Avellino', 0, 'Turris', 0, 1.45, '30/11/2023 09:13', 1.41, '03/12/2023 18:51', 3.94, '30/11/2023 09:13', 4.54, '03/12/2023 20:21', 6.72, '30/11/2023 09:13', 8, '03/12/2023 20:21', 'https://www.betexplorer.com/football/italy/serie-c-group-c/avellino-turris/CMpGl0PP/')
Set-RowFV 157 @('Brindisi', 0, 'Crotone', 2, 4.08, '30/11/2023 09:13', 3.64, '03/12/2023 20:38', 3.1, '30/11/2023 09:13', 3.41, '03/12/2023 20:38', 1.9, '30/11/2023 09:13', 2.04, '03/12/2023 20:38', 'https://www.betexplorer.com/football/italy/serie-c-group-c/brindisi-crotone/U1XspKXt/')
Set-RowFV 158 @('Potenza', 1, 'Taranto', 2, 2.39, '30/11/2023 09:13', 2.73, '04/12/2023 20:35', 2.81, '30/11/2023 09:13', 2.81, '04/12/2023 20:41', 3.14, '30/11/2023 09:13', 2.98, '04/12/2023 20:41', 'https://www.betexplorer.com/football/italy/serie-c-group-c/potenza-taranto/htI7vG1I/')
Set-RowFV 159 @('Casertana', 2, 'Foggia', 1, 2.09, '30/11/2023 09:13', 2.14, '04/12/2023 20:33', 2.97, '30/11/2023 09:13', 3.09, '04/12/2023 20:40', 3.59, '30/11/2023 09:13', 3.74, '04/12/2023 20:40', 'https://www.betexplorer.com/football/italy/serie-c-group-c/casertana-foggia/pbPPSTE2/')
Set-RowFV 160 @('Latina', 0, 'Sorrento', 2, 1.83, '30/11/2023 09:13', 1.68, '04/12/2023 20:40', 3.09, '30/11/2023 09:13', 3.33, '04/12/2023 20:41', 4.51, '30/11/2023 09:13', 6.2, '04/12/2023 20:40', 'https://www.betexplorer.com/football/italy/serie-c-group-c/latina-sorrento/4KPfsIHa/')
Set-RowFV 161 @('Monopoli', 1, 'Giugliano', 3, 1.92, '30/11/2023 09:13', 1.92, '04/12/2023 20:41', 3.15, '30/11/2023 09:13', 3.31, '04/12/2023 20:41', 3.95, '30/11/2023 09:13', 4.28, '04/12/2023 20:41', 'https://www.betexplorer.com/football/italy/serie-c-group-c/monopoli-giugliano/bsPbtxX5/')
Set-RowFV 165 @('Turris', 1, 'Latina', 1, 2.72, '08/12/2023 07:12', 2.75, '09/12/2023 20:08', 2.84, '08/12/2023 07:12', 3.01, '09/12/2023 20:12', 2.67, '08/12/2023 07:12', 2.77, '09/12/2023 20:12', 'https://www.betexplorer.com/football/italy/serie-c-group-c/turris-latina/IcKYhFGB/')
Set-RowFV 166 @('ACR Messina', 1, 'Catania', 0, 4.58, '08/12/2023 07:12', 5.27, '09/12/2023 20:36', 3.19, '08/12/2023 07:12', 3.4, '09/12/2023 20:36', 1.78, '08/12/2023 07:12', 1.75, '09/12/2023 20:36', 'https://www.betexplorer.com/football/italy/serie-c-group-c/acr-messina-catania/WWGBwzHO/')
Set-RowFV 178 @('Audace Cerignola', 2, 'Foggia', 0, 2.31, '14/12/2023 09:13', 2.57, '18/12/2023 20:36', 2.84, '14/12/2023 09:13', 3.08, '18/12/2023 20:36', 3.24, '14/12/2023 09:13', 2.91, '18/12/2023 20:36', 'https://www.betexplorer.com/football/italy/serie-c-group-c/audace-cerignola-foggia/lps6z70p/')
Set-RowFV 179 @('Casertana', 3, 'Giugliano', 1, 1.7, '14/12/2023 09:13', 1.89, '18/12/2023 20:43', 3.37, '14/12/2023 09:13', 3.79, '18/12/2023 20:43', 4.82, '14/12/2023 09:13', 3.78, '18/12/2023 20:43', 'https://www.betexplorer.com/football/italy/serie-c-group-c/casertana-giugliano/hOSXQ7qF/')
Set-RowFV 180 @('Catania', 0, 'Sorrento', 1, 1.63, '14/12/2023 09:13', 1.65, '18/12/2023 20:40', 3.33, '14/12/2023 09:13', 3.46, '18/12/2023 20:40', 5.64, '14/12/2023 09:13', 6.18, '18/12/2023 20:40', 'https://www.betexplorer.com/football/italy/serie-c-group-c/catania-sorrento/rqOdcU7c/')
Set-RowFV 181 @('Latina', 0, 'Benevento', 0, 2.91, '14/12/2023 09:13', 3.28, '18/12/2023 20:44', 2.79, '14/12/2023 09:13', 2.79, '18/12/2023 20:44', 2.56, '14/12/2023 09:13', 2.54, '18/12/2023 20:44', 'https://www.betexplorer.com/football/italy/serie-c-group-c/latina-benevento/j9R4e8xA/')

# --- Append the 5 new matches scraped for 22-12-2023 (rows 182-186). ---
Add-NewRow 182 181 45282.77083333334 @('Picerno', 0, 'Juve Stabia', 2, 2.39, '20/12/2023 18:42', 2.89, '22/12/2023 18:28', 2.73, '20/12/2023 18:42', 2.64, '22/12/2023 18:28', 3.23, '20/12/2023 18:42', 3.01, '22/12/2023 18:28', 'https://www.betexplorer.com/football/italy/serie-c-group-c/picerno-juve-stabia/nysq9R69/')
Add-NewRow 183 182 45282.77083333334 @('Turris', 3, 'Brindisi', 1, 2, '20/12/2023 18:42', 1.93, '22/12/2023 16:36', 3.05, '20/12/2023 18:42', 3.42, '22/12/2023 18:02', 3.79, '20/12/2023 18:42', 4.02, '22/12/2023 18:02', 'https://www.betexplorer.com/football/italy/serie-c-group-c/turris-brindisi/Yiui75yM/')
Add-NewRow 184 183 45282.86458333334 @('ACR Messina', 1, 'Monopoli', 1, 2.65, '20/12/2023 21:12', 2.17, '22/12/2023 20:44', 2.74, '20/12/2023 21:12', 2.98, '22/12/2023 20:18', 2.86, '20/12/2023 21:12', 3.82, '22/12/2023 20:44', 'https://www.betexplorer.com/football/italy/serie-c-group-c/acr-messina-monopoli/tv6qlQUq/')
Add-NewRow 185 184 45282.86458333334 @('Crotone', 0, 'Avellino', 1, 2.51, '20/12/2023 21:12', 2.89, '22/12/2023 20:41', 2.79, '20/12/2023 21:12', 3.03, '22/12/2023 20:39', 2.96, '20/12/2023 21:12', 2.62, '22/12/2023 20:41', 'https://www.betexplorer.com/football/italy/serie-c-group-c/crotone-avellino/x6ahn4ad/')
Add-NewRow 186 185 45282.86458333334 @('Taranto', 2, 'Latina', 1, 2.16, '20/12/2023 21:12', 1.93, '22/12/2023 20:40', 2.83, '20/12/2023 21:12', 2.82, '22/12/2023 20:40', 3.61, '20/12/2023 21:12', 5.43, '22/12/2023 20:40', 'https://www.betexplorer.com/football/italy/serie-c-group-c/taranto-latina/QRum8oMF/')
